$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 22. This shifts the existing rows
# 22:94 down to 23:95 (carrying their values/styles with them), which
# matches the target diff exactly (every row N>=23 ends up holding what
# used to be row N-1's data, and a new last row 95 appears holding the
# old row 94's data).
$ws.Rows("22").Insert()

# Populate the newly inserted row 22 with the new weekly price-report
# entry. Columns A,B,C,E,F,G,H,I,J,K are constant across the whole
# dataset; L (Calidad) and R (Origen) keep the same values this series
# already had for this slot ("Primera" / "Provincia de Limar\u00ed").
$ws.Cells.Item(22, 1).Value = 3
$ws.Cells.Item(22, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(22, 3).Value = "Coquimbo"
$ws.Cells.Item(22, 4).Value = 44592
$ws.Cells.Item(22, 5).Value = 5
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100107
$ws.Cells.Item(22, 8).Value = "Otros"
$ws.Cells.Item(22, 9).Value = 100107011
$ws.Cells.Item(22, 10).Value = "Tuna"
$ws.Cells.Item(22, 11).Value = "Sin especificar"
$ws.Cells.Item(22, 12).Value = "Primera"
$ws.Cells.Item(22, 13).Value = 50
$ws.Cells.Item(22, 14).Value = 20000
$ws.Cells.Item(22, 15).Value = 20000
$ws.Cells.Item(22, 16).Value = 20000
$ws.Cells.Item(22, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(22, 18).Value = "Provincia de Limar" + [char]0x00ED
$ws.Cells.Item(22, 19).Value = 1000
$ws.Cells.Item(22, 20).Value = 20
